$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Rows.Item(1).Delete()
$ws2.Rows.Item(1).Delete()
